# Update the recommended ticket price from $54/$56 to $64 across the deck.

$p = $ppt.ActivePresentation

# Slide 3 ("Answer:") - Content Placeholder 2:
#   "Raise ticket prices to $54." -> "Raise ticket prices to $64."
# Split into three runs, matching how PowerPoint records an in-place edit
# of a run's text (InsertAfter creates new sibling runs).
$slide3 = $p.Slides.Item(3)
$tr3 = $slide3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "Raise ticket prices to "
$run3b = $tr3.InsertAfter("$64")
$run3c = $run3b.InsertAfter(".")

# Slide 4 ("Modelling") - Content Placeholder 2, third paragraph:
#   "Prices for all three were similar - ~$56 recommended."
#   -> "Prices for all three were similar - ~$64 recommended."
$slide4 = $p.Slides.Item(4)
$tr4 = $slide4.Shapes.Item(2).TextFrame.TextRange
$para4 = $tr4.Paragraphs(3)
$para4.Text = "Prices for all three were similar - "
$run4b = $para4.InsertAfter("~$64 ")
$run4c = $run4b.InsertAfter("recommended.")
